$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 04.02.2022 13:00"

# Update row 3 (Tesco): Delta Cena and Old Datum become real numeric values
$ws.Range("D3").Value = 0.01
$ws.Range("E3").Value = 44596.53134259259
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
